$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new values are plain strings (non-ambiguous as text)
$plainUpdates = @{
    "D2" = "61.453.55"
    "E2" = "  -3.03%  "
    "D3" = "2.986.31"
    "E3" = "  -2.92%  "
    "E4" = "  -0.09%  "
    "E5" = "  +0.97%  "
    "E6" = "  -5.62%  "
    "E7" = "  -0.15%  "
    "D8" = "2.979.94"
    "E8" = "  -2.90%  "
    "E9" = "  -1.57%  "
    "E10" = "  -6.39%  "
    "E11" = "  -8.08%  "
    "E12" = "  -2.70%  "
    "E13" = "  -2.93%  "
    "E14" = "  -2.90%  "
    "D15" = "3.465.73"
    "E15" = "  -3.17%  "
    "D16" = "61.569.15"
    "E16" = "  -2.93%  "
    "E17" = "  -2.69%  "
    "D18" = "2.984.93"
    "E18" = "  -3.14%  "
    "E19" = "  -0.43%  "
    "E20" = "  +0.94%  "
    "E21" = "  -2.83%  "
    "E22" = "  -5.21%  "
    "E23" = "  -1.51%  "
    "E24" = "  +2.26%  "
    "E25" = "  -1.89%  "
    "E26" = "  -0.11%  "
    "E27" = "  -0.42%  "
    "E28" = "  -3.91%  "
    "E29" = "  -0.14%  "
    "E30" = "  +0.09%  "
    "E31" = "  -2.84%  "
    "E32" = "  -2.83%  "
    "E33" = "  -0.36%  "
    "E34" = "  -0.34%  "
    "E35" = "  -7.56%  "
    "E36" = "  -2.57%  "
    "E37" = "  -7.94%  "
    "D38" = "3.114.04"
    "E38" = "  -4.56%  "
    "E39" = "  -0.99%  "
    "E40" = "  -6.13%  "
    "E41" = "  -1.74%  "
    "E42" = "  -0.99%  "
    "E43" = "  -0.04%  "
    "E44" = "  -10.26%  "
    "E45" = "  -0.79%  "
    "E46" = "  -5.13%  "
    "E47" = "  -1.61%  "
    "E48" = "  -4.22%  "
    "E49" = "  +9.35%  "
    "E50" = "  -7.69%  "
    "E51" = "  -9.59%  "
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Cells whose new values look numeric but must remain text (e.g. "0.999", "130.68")
# Prefix with an apostrophe to force text entry, then reset the style so no
# residual number-format/quote-prefix styling is left on the cell.
$textUpdates = @{
    "D5" = "549.64"
    "D6" = "130.68"
    "D9" = "0.490"
    "D11" = "0.144"
    "D12" = "0.443"
    "D14" = "33.84"
    "D19" = "6.59"
    "D20" = "478.88"
    "D23" = "6.97"
    "D24" = "80.31"
    "D25" = "11.98"
    "D26" = "0.999"
    "D27" = "2.71"
    "D28" = "7.64"
    "D29" = "0.999"
    "D31" = "25.44"
    "D33" = "2.30"
    "D34" = "5.47"
    "D35" = "54.49"
    "D37" = "448.47"
    "D39" = "0.0788"
    "D42" = "8.05"
    "D44" = "2.32"
    "D48" = "1.94"
    "D50" = "113.77"
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = "'" + $textUpdates[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}

